# Generate Report for Handback
# Updates the "Latest HO Xliff Generate Date" / "Correspond Handoff Datetime" /
# "Correspond Handback DateTime" timestamps for the 66255b7e... file across
# the Overview, zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

# Overview sheet: "Latest HO Xliff Generate Date" for row 2 (66255b7e...)
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-09-03 13:10:15"

# zh-cn sheet: Correspond Handoff Datetime / Correspond Handback DateTime for row 2
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-09-03 13:10:06"
$wsZhCn.Range("K2").Value = "2016-09-03 13:10:27"

# de-de sheet: Correspond Handoff Datetime / Correspond Handback DateTime for row 2
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H2").Value = "2016-09-03 13:10:15"
$wsDeDe.Range("K2").Value = "2016-09-03 13:10:34"
